$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.385.73"
$ws.Range("E2").Value = "  +0.27%  "

$ws.Range("D3").Value = "1.868.49"
$ws.Range("E3").Value = "  -0.33%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "243.79"
$ws.Range("E5").Value = "  +0.54%  "

$ws.Range("D6").Value = "0.7053"
$ws.Range("E6").Value = "  -0.99%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "0.07923"
$ws.Range("E8").Value = "  -1.08%  "

$ws.Range("D9").Value = "0.3134"
$ws.Range("E9").Value = "  -0.62%  "

$ws.Range("D10").Value = "24.47"
$ws.Range("E10").Value = "  -2.00%  "

$ws.Range("D11").Value = "0.07837"
$ws.Range("E11").Value = "  -4.71%  "

$ws.Range("D12").Value = "1.863.74"
$ws.Range("E12").Value = "  -1.12%  "

$ws.Range("D13").Value = "93.75"
$ws.Range("E13").Value = "  -1.08%  "

$ws.Range("D14").Value = "5.182"
$ws.Range("E14").Value = "  -1.23%  "

$ws.Range("D15").Value = "0.7014"
$ws.Range("E15").Value = "  -1.40%  "

$ws.Range("D16").Value = "6.524"
$ws.Range("E16").Value = "  +1.69%  "

$ws.Range("D17").Value = "0.000008395"
$ws.Range("E17").Value = "  -1.79%  "

$ws.Range("D18").Value = "29.386.80"
$ws.Range("E18").Value = "  +0.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.90"
$ws.Range("E19").Value = "  +4.21%  "

$ws.Range("D20").Value = "2.121.40"
$ws.Range("E20").Value = "  -1.15%  "

$ws.Range("D21").Value = "13.11"
$ws.Range("E21").Value = "  -0.98%  "

$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("D23").Value = "7.652"
$ws.Range("E23").Value = "  -1.51%  "

$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").Value = "0.1556"
$ws.Range("E25").Value = "  -0.27%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.010"
$ws.Range("E26").Value = "  -0.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.90"

$ws.Range("D28").Value = "18.83"
$ws.Range("E28").Value = "  +1.69%  "

$ws.Range("D29").Value = "1.502"
$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.320"
$ws.Range("E30").Value = "  -2.08%  "

$ws.Range("D31").Value = "4.252"
$ws.Range("E31").Value = "  -1.26%  "

$ws.Range("D32").Value = "1.212"
$ws.Range("E32").Value = "  +3.90%  "

$ws.Range("D33").Value = "0.05297"
$ws.Range("E33").Value = "  -1.47%  "

$ws.Range("D34").Value = "1.897"
$ws.Range("E34").Value = "  -2.21%  "

$ws.Range("D35").Value = "0.7496"
$ws.Range("E35").Value = "  -1.78%  "

$ws.Range("D36").Value = "1.174"
$ws.Range("E36").Value = "  -0.26%  "

$ws.Range("E37").Value = "  +0.86%  "

$ws.Range("D38").Value = "1.290.72"

$ws.Range("D39").Value = "0.01879"
$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("D40").Value = "2.764"
$ws.Range("E40").Value = "  +0.43%  "

$ws.Range("E41").Value = "  -2.20%  "

$ws.Range("D44").Value = "71.07"
$ws.Range("E44").Value = "  -4.01%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("E46").Value = "  -3.41%  "

$ws.Range("D47").Value = "2.021.20"

$ws.Range("D48").Value = "1.799"
$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("E49").Value = "  +1.26%  "

$ws.Range("D50").Value = "0.5179"
$ws.Range("E50").Value = "  -0.82%  "

$ws.Range("D51").Value = "0.4303"
$ws.Range("E51").Value = "  -1.07%  "

# Row 42/43: FraxShare and Quant swap positions with updated values
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "108.82"
$ws.Range("E42").Value = "  -3.51%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "6.011"
$ws.Range("E43").Value = "  -7.20%  "

